$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999817537866
$ws.Range("A2").Value = 0.99835703773066364
$ws.Range("A3").Value = 0.99312433938167111
$ws.Range("A4").Value = 0.99538545395598477
$ws.Range("A5").Value = 0.98446652423590386
$ws.Range("A6").Value = 0.96019022919848662
$ws.Range("A7").Value = 0.95437912392683732
$ws.Range("A8").Value = 0.945375502941462
$ws.Range("A9").Value = 0.9343487815133481
$ws.Range("A10").Value = 0.92423224848715513
$ws.Range("A11").Value = 0.92263588361826687
$ws.Range("A12").Value = 0.9198557521521471
$ws.Range("A13").Value = 0.90856761187951363
$ws.Range("A14").Value = 0.90439993205005709
$ws.Range("A15").Value = 0.90180825626434558
$ws.Range("A16").Value = 0.89930164093596832
$ws.Range("A17").Value = 0.89559355467595658
$ws.Range("A18").Value = 0.89448461029216975
$ws.Range("A19").Value = 0.99730376921226416
$ws.Range("A20").Value = 0.99018656759115165
$ws.Range("A21").Value = 0.98878804487214622
$ws.Range("A22").Value = 0.98752353483799671
$ws.Range("A23").Value = 0.98609005334479982
$ws.Range("A24").Value = 0.97306971363661998
$ws.Range("A25").Value = 0.96661279528473565
$ws.Range("A26").Value = 0.95820205602409825
$ws.Range("A27").Value = 0.95568025116968391
$ws.Range("A28").Value = 0.94728316853121453
$ws.Range("A29").Value = 0.94172371925148757
$ws.Range("A30").Value = 0.93712790423649517
$ws.Range("A31").Value = 0.93755417376569683
$ws.Range("A32").Value = 0.93587483313782505
$ws.Range("A33").Value = 0.9353548050518814
